$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$q3.Range("B1:H1").NumberFormat = "@"
$q3.Range("B2:G26").NumberFormat = "@"
$q3.Range("B27:F27").NumberFormat = "@"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "000991"
$q3.Range("C2").Value = "工银战略转型股票A"
$q3.Range("D2").Value = "44.53"
$q3.Range("E2").Value = "93.90"
$q3.Range("F2").Value = "4.40"
$q3.Range("G2").Value = "1.9593"
$q3.Range("H2").Value = 7
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "009428"
$q3.Range("C3").Value = "鹏扬景沣六个月持有期混合A"
$q3.Range("D3").Value = "35.59"
$q3.Range("E3").Value = "24.38"
$q3.Range("F3").Value = "1.85"
$q3.Range("G3").Value = "0.6584"
$q3.Range("H3").Value = 4
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "000971"
$q3.Range("C4").Value = "诺安新经济股票"
$q3.Range("D4").Value = "14.81"
$q3.Range("E4").Value = "87.07"
$q3.Range("F4").Value = "3.35"
$q3.Range("G4").Value = "0.4961"
$q3.Range("H4").Value = 9
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "002670"
$q3.Range("C5").Value = "万家沪深300指数增强A"
$q3.Range("D5").Value = "20.85"
$q3.Range("E5").Value = "94.06"
$q3.Range("F5").Value = "2.30"
$q3.Range("G5").Value = "0.4796"
$q3.Range("H5").Value = 4
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "011818"
$q3.Range("C6").Value = "鹏扬景阳一年持有期混合A"
$q3.Range("D6").Value = "22.73"
$q3.Range("E6").Value = "23.35"
$q3.Range("F6").Value = "1.76"
$q3.Range("G6").Value = "0.4000"
$q3.Range("H6").Value = 4
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "002671"
$q3.Range("C7").Value = "万家沪深300指数增强C"
$q3.Range("D7").Value = "10.38"
$q3.Range("E7").Value = "94.06"
$q3.Range("F7").Value = "2.30"
$q3.Range("G7").Value = "0.2387"
$q3.Range("H7").Value = 4
$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "011521"
$q3.Range("C8").Value = "鹏扬景源一年持有期混合A"
$q3.Range("D8").Value = "13.87"
$q3.Range("E8").Value = "23.83"
$q3.Range("F8").Value = "1.53"
$q3.Range("G8").Value = "0.2122"
$q3.Range("H8").Value = 4
$q3.Range("A9").Value = 7
$q3.Range("B9").Value = "009064"
$q3.Range("C9").Value = "鹏扬景沃六个月持有期混合A"
$q3.Range("D9").Value = "12.47"
$q3.Range("E9").Value = "25.11"
$q3.Range("F9").Value = "1.23"
$q3.Range("G9").Value = "0.1534"
$q3.Range("H9").Value = 5
$q3.Range("A10").Value = 8
$q3.Range("B10").Value = "011473"
$q3.Range("C10").Value = "工银战略转型股票C"
$q3.Range("D10").Value = "2.83"
$q3.Range("E10").Value = "93.90"
$q3.Range("F10").Value = "4.40"
$q3.Range("G10").Value = "0.1245"
$q3.Range("H10").Value = 7
$q3.Range("A11").Value = 9
$q3.Range("B11").Value = "393001"
$q3.Range("C11").Value = "中海优势精选灵活配置混合"
$q3.Range("D11").Value = "1.50"
$q3.Range("E11").Value = "78.67"
$q3.Range("F11").Value = "7.85"
$q3.Range("G11").Value = "0.1178"
$q3.Range("H11").Value = 4
$q3.Range("A12").Value = 10
$q3.Range("B12").Value = "009429"
$q3.Range("C12").Value = "鹏扬景沣六个月持有期混合C"
$q3.Range("D12").Value = "5.55"
$q3.Range("E12").Value = "24.38"
$q3.Range("F12").Value = "1.85"
$q3.Range("G12").Value = "0.1027"
$q3.Range("H12").Value = 4
$q3.Range("A13").Value = 11
$q3.Range("B13").Value = "009065"
$q3.Range("C13").Value = "鹏扬景沃六个月持有期混合C"
$q3.Range("D13").Value = "5.83"
$q3.Range("E13").Value = "25.11"
$q3.Range("F13").Value = "1.23"
$q3.Range("G13").Value = "0.0717"
$q3.Range("H13").Value = 5
$q3.Range("A14").Value = 12
$q3.Range("B14").Value = "009130"
$q3.Range("C14").Value = "鹏扬景恒六个月持有期混合A"
$q3.Range("D14").Value = "4.75"
$q3.Range("E14").Value = "24.46"
$q3.Range("F14").Value = "1.35"
$q3.Range("G14").Value = "0.0641"
$q3.Range("H14").Value = 5
$q3.Range("A15").Value = 13
$q3.Range("B15").Value = "002849"
$q3.Range("C15").Value = "金信智能中国2025灵活配置混合"
$q3.Range("D15").Value = "1.09"
$q3.Range("E15").Value = "78.83"
$q3.Range("F15").Value = "4.93"
$q3.Range("G15").Value = "0.0537"
$q3.Range("H15").Value = 6
$q3.Range("A16").Value = 14
$q3.Range("B16").Value = "012708"
$q3.Range("C16").Value = "东方红中证东方红红利低波动指数A"
$q3.Range("D16").Value = "3.27"
$q3.Range("E16").Value = "93.80"
$q3.Range("F16").Value = "1.56"
$q3.Range("G16").Value = "0.0510"
$q3.Range("H16").Value = 10
$q3.Range("A17").Value = 15
$q3.Range("B17").Value = "011522"
$q3.Range("C17").Value = "鹏扬景源一年持有期混合C"
$q3.Range("D17").Value = "2.75"
$q3.Range("E17").Value = "23.83"
$q3.Range("F17").Value = "1.53"
$q3.Range("G17").Value = "0.0421"
$q3.Range("H17").Value = 4
$q3.Range("A18").Value = 16
$q3.Range("B18").Value = "001780"
$q3.Range("C18").Value = "诺安改革趋势灵活配置混合"
$q3.Range("D18").Value = "0.54"
$q3.Range("E18").Value = "85.03"
$q3.Range("F18").Value = "6.10"
$q3.Range("G18").Value = "0.0329"
$q3.Range("H18").Value = 6
$q3.Range("A19").Value = 17
$q3.Range("B19").Value = "515300"
$q3.Range("C19").Value = "嘉实沪深300红利低波动ETF"
$q3.Range("D19").Value = "0.94"
$q3.Range("E19").Value = "99.19"
$q3.Range("F19").Value = "2.91"
$q3.Range("G19").Value = "0.0274"
$q3.Range("H19").Value = 9
$q3.Range("A20").Value = 18
$q3.Range("B20").Value = "011819"
$q3.Range("C20").Value = "鹏扬景阳一年持有期混合C"
$q3.Range("D20").Value = "1.34"
$q3.Range("E20").Value = "23.35"
$q3.Range("F20").Value = "1.76"
$q3.Range("G20").Value = "0.0236"
$q3.Range("H20").Value = 4
$q3.Range("A21").Value = 19
$q3.Range("B21").Value = "009131"
$q3.Range("C21").Value = "鹏扬景恒六个月持有期混合C"
$q3.Range("D21").Value = "1.54"
$q3.Range("E21").Value = "24.46"
$q3.Range("F21").Value = "1.35"
$q3.Range("G21").Value = "0.0208"
$q3.Range("H21").Value = 5
$q3.Range("A22").Value = 20
$q3.Range("B22").Value = "510060"
$q3.Range("C22").Value = "工银上证央企50ETF"
$q3.Range("D22").Value = "0.76"
$q3.Range("E22").Value = "97.77"
$q3.Range("F22").Value = "2.63"
$q3.Range("G22").Value = "0.0200"
$q3.Range("H22").Value = 9
$q3.Range("A23").Value = 21
$q3.Range("B23").Value = "012709"
$q3.Range("C23").Value = "东方红中证东方红红利低波动指数C"
$q3.Range("D23").Value = "0.67"
$q3.Range("E23").Value = "93.80"
$q3.Range("F23").Value = "1.56"
$q3.Range("G23").Value = "0.0105"
$q3.Range("H23").Value = 10
$q3.Range("A24").Value = 22
$q3.Range("B24").Value = "007751"
$q3.Range("C24").Value = "景顺长城中证沪港深红利成长低波动指数A"
$q3.Range("D24").Value = "0.67"
$q3.Range("E24").Value = "90.27"
$q3.Range("F24").Value = "1.51"
$q3.Range("G24").Value = "0.0101"
$q3.Range("H24").Value = 10
$q3.Range("A25").Value = 23
$q3.Range("B25").Value = "009384"
$q3.Range("C25").Value = "摩根士丹利华鑫MSCI中国A股指数增强A"
$q3.Range("D25").Value = "0.39"
$q3.Range("E25").Value = "90.98"
$q3.Range("F25").Value = "1.30"
$q3.Range("G25").Value = "0.0051"
$q3.Range("H25").Value = 3
$q3.Range("A26").Value = 24
$q3.Range("B26").Value = "007760"
$q3.Range("C26").Value = "景顺长城中证沪港深红利成长低波动指数C"
$q3.Range("D26").Value = "0.06"
$q3.Range("E26").Value = "90.27"
$q3.Range("F26").Value = "1.51"
$q3.Range("G26").Value = "0.0009"
$q3.Range("H26").Value = 10
$q3.Range("A27").Value = 25
$q3.Range("B27").Value = "014866"
$q3.Range("C27").Value = "摩根士丹利华鑫MSCI中国A股指数增强C"
$q3.Range("D27").Value = "0.00"
$q3.Range("E27").Value = "90.98"
$q3.Range("F27").Value = "1.30"
$q3.Range("G27").Value = 0
$q3.Range("H27").Value = 3


# Insert new row for 2022-Q3 at the top of the "总计" (Total) summary table,
# shifting all the existing quarter rows down by one.
$total.Rows.Item(2).Insert()
# Copy cell A3's style (bold/centered/bordered index style) onto the newly
# inserted A2 so the row-index column keeps its original formatting.
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 26
$total.Range("D2").Value = 5.38
